$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Sdc2"
$ws.Range("C2").Value = "Ptprj"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = [double]"3"
$ws.Range("F2").Value = [double]"1"
$ws.Range("G2").Value = [double]"3.057109"
$ws.Range("H2").Value = [double]"9.171327"
$ws.Range("I2").Value = [double]"0.02694952608666365"
$ws.Range("J2").Value = [double]"0.02694952608666365"
$ws.Range("K2").Value = [double]"3"
$ws.Range("L2").Value = [double]"1"
$ws.Range("M2").Value = [double]"1.815493333333333"
$ws.Range("N2").Value = [double]"5.44648"
$ws.Range("O2").Value = [double]"0.02449420924905278"
$ws.Range("P2").Value = [double]"0.02449420924905277"
$ws.Range("Q2").Value = [double]"5.550161008773333"
$ws.Range("R2").Value = [double]"49.95144907896"
$ws.Range("S2").Value = [double]"0.0006601073311295459"
$ws.Range("T2").Value = [double]"0.0006601073311295458"

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Sdc2"
$ws.Range("C3").Value = "Ptprj"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = [double]"3"
$ws.Range("F3").Value = [double]"1"
$ws.Range("G3").Value = [double]"3.057109"
$ws.Range("H3").Value = [double]"9.171327"
$ws.Range("I3").Value = [double]"0.02694952608666365"
$ws.Range("J3").Value = [double]"0.02694952608666365"
$ws.Range("K3").Value = [double]"3"
$ws.Range("L3").Value = [double]"1"
$ws.Range("M3").Value = [double]"5.684019666666667"
$ws.Range("N3").Value = [double]"17.052059"
$ws.Range("O3").Value = [double]"0.07668745708663094"
$ws.Range("P3").Value = [double]"0.07668745708663093"
$ws.Range("Q3").Value = [double]"17.37666767914367"
$ws.Range("R3").Value = [double]"156.390009112293"
$ws.Range("S3").Value = [double]"0.00206669062527606"
$ws.Range("T3").Value = [double]"0.002066690625276059"

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Sdc2"
$ws.Range("C4").Value = "Ptprj"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = [double]"3"
$ws.Range("F4").Value = [double]"1"
$ws.Range("G4").Value = [double]"3.057109"
$ws.Range("H4").Value = [double]"9.171327"
$ws.Range("I4").Value = [double]"0.02694952608666365"
$ws.Range("J4").Value = [double]"0.02694952608666365"
$ws.Range("K4").Value = [double]"3"
$ws.Range("L4").Value = [double]"1"
$ws.Range("M4").Value = [double]"61.37607633333332"
$ws.Range("N4").Value = [double]"184.128229"
$ws.Range("O4").Value = [double]"0.8280715929891429"
$ws.Range("P4").Value = [double]"0.8280715929891428"
$ws.Range("Q4").Value = [double]"187.6333553433203"
$ws.Range("R4").Value = [double]"1688.700198089883"
$ws.Range("S4").Value = [double]"0.02231613699688603"
$ws.Range("T4").Value = [double]"0.02231613699688603"

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Sdc2"
$ws.Range("C5").Value = "Ptprj"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = [double]"3"
$ws.Range("F5").Value = [double]"1"
$ws.Range("G5").Value = [double]"3.057109"
$ws.Range("H5").Value = [double]"9.171327"
$ws.Range("I5").Value = [double]"0.02694952608666365"
$ws.Range("J5").Value = [double]"0.02694952608666365"
$ws.Range("K5").Value = [double]"3"
$ws.Range("L5").Value = [double]"1"
$ws.Range("M5").Value = [double]"5.243698"
$ws.Range("N5").Value = [double]"15.731094"
$ws.Range("O5").Value = [double]"0.07074674067517345"
$ws.Range("P5").Value = [double]"0.07074674067517343"
$ws.Range("Q5").Value = [double]"16.030556349082"
$ws.Range("R5").Value = [double]"144.275007141738"
$ws.Range("S5").Value = [double]"0.001906591133372015"
$ws.Range("T5").Value = [double]"0.001906591133372015"

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Sdc2"
$ws.Range("C6").Value = "Ptprj"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = [double]"3"
$ws.Range("F6").Value = [double]"1"
$ws.Range("G6").Value = [double]"89.02756"
$ws.Range("H6").Value = [double]"267.08268"
$ws.Range("I6").Value = [double]"0.7848102735793893"
$ws.Range("J6").Value = [double]"0.7848102735793893"
$ws.Range("K6").Value = [double]"3"
$ws.Range("L6").Value = [double]"1"
$ws.Range("M6").Value = [double]"1.815493333333333"
$ws.Range("N6").Value = [double]"5.44648"
$ws.Range("O6").Value = [double]"0.02449420924905278"
$ws.Range("P6").Value = [double]"0.02449420924905277"
$ws.Range("Q6").Value = [double]"161.6289416629333"
$ws.Range("R6").Value = [double]"1454.6604749664"
$ws.Range("S6").Value = [double]"0.01922330706185992"
$ws.Range("T6").Value = [double]"0.01922330706185992"

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Sdc2"
$ws.Range("C7").Value = "Ptprj"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = [double]"3"
$ws.Range("F7").Value = [double]"1"
$ws.Range("G7").Value = [double]"89.02756"
$ws.Range("H7").Value = [double]"267.08268"
$ws.Range("I7").Value = [double]"0.7848102735793893"
$ws.Range("J7").Value = [double]"0.7848102735793893"
$ws.Range("K7").Value = [double]"3"
$ws.Range("L7").Value = [double]"1"
$ws.Range("M7").Value = [double]"5.684019666666667"
$ws.Range("N7").Value = [double]"17.052059"
$ws.Range("O7").Value = [double]"0.07668745708663094"
$ws.Range("P7").Value = [double]"0.07668745708663093"
$ws.Range("Q7").Value = [double]"506.0344019153467"
$ws.Range("R7").Value = [double]"4554.309617238119"
$ws.Range("S7").Value = [double]"0.0601851041762665"
$ws.Range("T7").Value = [double]"0.06018510417626649"

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Sdc2"
$ws.Range("C8").Value = "Ptprj"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = [double]"3"
$ws.Range("F8").Value = [double]"1"
$ws.Range("G8").Value = [double]"89.02756"
$ws.Range("H8").Value = [double]"267.08268"
$ws.Range("I8").Value = [double]"0.7848102735793893"
$ws.Range("J8").Value = [double]"0.7848102735793893"
$ws.Range("K8").Value = [double]"3"
$ws.Range("L8").Value = [double]"1"
$ws.Range("M8").Value = [double]"61.37607633333332"
$ws.Range("N8").Value = [double]"184.128229"
$ws.Range("O8").Value = [double]"0.8280715929891429"
$ws.Range("P8").Value = [double]"0.8280715929891428"
$ws.Range("Q8").Value = [double]"5464.162318330412"
$ws.Range("R8").Value = [double]"49177.46086497371"
$ws.Range("S8").Value = [double]"0.6498790934371299"
$ws.Range("T8").Value = [double]"0.6498790934371298"

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Sdc2"
$ws.Range("C9").Value = "Ptprj"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = [double]"3"
$ws.Range("F9").Value = [double]"1"
$ws.Range("G9").Value = [double]"89.02756"
$ws.Range("H9").Value = [double]"267.08268"
$ws.Range("I9").Value = [double]"0.7848102735793893"
$ws.Range("J9").Value = [double]"0.7848102735793893"
$ws.Range("K9").Value = [double]"3"
$ws.Range("L9").Value = [double]"1"
$ws.Range("M9").Value = [double]"5.243698"
$ws.Range("N9").Value = [double]"15.731094"
$ws.Range("O9").Value = [double]"0.07074674067517345"
$ws.Range("P9").Value = [double]"0.07074674067517343"
$ws.Range("Q9").Value = [double]"466.83363831688"
$ws.Range("R9").Value = [double]"4201.502744851919"
$ws.Range("S9").Value = [double]"0.05552276890413298"
$ws.Range("T9").Value = [double]"0.05552276890413297"

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Sdc2"
$ws.Range("C10").Value = "Ptprj"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = [double]"2"
$ws.Range("F10").Value = [double]"0.6666666666666666"
$ws.Range("G10").Value = [double]"0.184005"
$ws.Range("H10").Value = [double]"0.5520149999999999"
$ws.Range("I10").Value = [double]"0.0016220709001794"
$ws.Range("J10").Value = [double]"0.0016220709001794"
$ws.Range("K10").Value = [double]"3"
$ws.Range("L10").Value = [double]"1"
$ws.Range("M10").Value = [double]"1.815493333333333"
$ws.Range("N10").Value = [double]"5.44648"
$ws.Range("O10").Value = [double]"0.02449420924905278"
$ws.Range("P10").Value = [double]"0.02449420924905277"
$ws.Range("Q10").Value = [double]"0.3340598507999999"
$ws.Range("R10").Value = [double]"3.0065386572"
$ws.Range("S10").Value = [double]"3.973134404579362e-05"
$ws.Range("T10").Value = [double]"3.973134404579361e-05"

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Sdc2"
$ws.Range("C11").Value = "Ptprj"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = [double]"2"
$ws.Range("F11").Value = [double]"0.6666666666666666"
$ws.Range("G11").Value = [double]"0.184005"
$ws.Range("H11").Value = [double]"0.5520149999999999"
$ws.Range("I11").Value = [double]"0.0016220709001794"
$ws.Range("J11").Value = [double]"0.0016220709001794"
$ws.Range("K11").Value = [double]"3"
$ws.Range("L11").Value = [double]"1"
$ws.Range("M11").Value = [double]"5.684019666666667"
$ws.Range("N11").Value = [double]"17.052059"
$ws.Range("O11").Value = [double]"0.07668745708663094"
$ws.Range("P11").Value = [double]"0.07668745708663093"
$ws.Range("Q11").Value = [double]"1.045888038765"
$ws.Range("R11").Value = [double]"9.412992348884998"
$ws.Range("S11").Value = [double]"0.0001243924925489805"
$ws.Range("T11").Value = [double]"0.0001243924925489805"

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Sdc2"
$ws.Range("C12").Value = "Ptprj"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = [double]"2"
$ws.Range("F12").Value = [double]"0.6666666666666666"
$ws.Range("G12").Value = [double]"0.184005"
$ws.Range("H12").Value = [double]"0.5520149999999999"
$ws.Range("I12").Value = [double]"0.0016220709001794"
$ws.Range("J12").Value = [double]"0.0016220709001794"
$ws.Range("K12").Value = [double]"3"
$ws.Range("L12").Value = [double]"1"
$ws.Range("M12").Value = [double]"61.37607633333332"
$ws.Range("N12").Value = [double]"184.128229"
$ws.Range("O12").Value = [double]"0.8280715929891429"
$ws.Range("P12").Value = [double]"0.8280715929891428"
$ws.Range("Q12").Value = [double]"11.293504925715"
$ws.Range("R12").Value = [double]"101.641544331435"
$ws.Range("S12").Value = [double]"0.001343190834252889"
$ws.Range("T12").Value = [double]"0.001343190834252888"

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Sdc2"
$ws.Range("C13").Value = "Ptprj"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = [double]"2"
$ws.Range("F13").Value = [double]"0.6666666666666666"
$ws.Range("G13").Value = [double]"0.184005"
$ws.Range("H13").Value = [double]"0.5520149999999999"
$ws.Range("I13").Value = [double]"0.0016220709001794"
$ws.Range("J13").Value = [double]"0.0016220709001794"
$ws.Range("K13").Value = [double]"3"
$ws.Range("L13").Value = [double]"1"
$ws.Range("M13").Value = [double]"5.243698"
$ws.Range("N13").Value = [double]"15.731094"
$ws.Range("O13").Value = [double]"0.07074674067517345"
$ws.Range("P13").Value = [double]"0.07074674067517343"
$ws.Range("Q13").Value = [double]"0.9648666504899999"
$ws.Range("R13").Value = [double]"8.68379985441"
$ws.Range("S13").Value = [double]"0.0001147562293317372"
$ws.Range("T13").Value = [double]"0.0001147562293317371"

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Sdc2"
$ws.Range("C14").Value = "Ptprj"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = [double]"3"
$ws.Range("F14").Value = [double]"1"
$ws.Range("G14").Value = [double]"21.16964733333333"
$ws.Range("H14").Value = [double]"63.508942"
$ws.Range("I14").Value = [double]"0.1866181294337677"
$ws.Range("J14").Value = [double]"0.1866181294337677"
$ws.Range("K14").Value = [double]"3"
$ws.Range("L14").Value = [double]"1"
$ws.Range("M14").Value = [double]"1.815493333333333"
$ws.Range("N14").Value = [double]"5.44648"
$ws.Range("O14").Value = [double]"0.02449420924905278"
$ws.Range("P14").Value = [double]"0.02449420924905277"
$ws.Range("Q14").Value = [double]"38.43335360268445"
$ws.Range("R14").Value = [double]"345.90018242416"
$ws.Range("S14").Value = [double]"0.004571063512017522"
$ws.Range("T14").Value = [double]"0.004571063512017521"

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Sdc2"
$ws.Range("C15").Value = "Ptprj"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = [double]"3"
$ws.Range("F15").Value = [double]"1"
$ws.Range("G15").Value = [double]"21.16964733333333"
$ws.Range("H15").Value = [double]"63.508942"
$ws.Range("I15").Value = [double]"0.1866181294337677"
$ws.Range("J15").Value = [double]"0.1866181294337677"
$ws.Range("K15").Value = [double]"3"
$ws.Range("L15").Value = [double]"1"
$ws.Range("M15").Value = [double]"5.684019666666667"
$ws.Range("N15").Value = [double]"17.052059"
$ws.Range("O15").Value = [double]"0.07668745708663094"
$ws.Range("P15").Value = [double]"0.07668745708663093"
$ws.Range("Q15").Value = [double]"120.3286917790642"
$ws.Range("R15").Value = [double]"1082.958226011578"
$ws.Range("S15").Value = [double]"0.0143112697925394"
$ws.Range("T15").Value = [double]"0.0143112697925394"

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Sdc2"
$ws.Range("C16").Value = "Ptprj"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = [double]"3"
$ws.Range("F16").Value = [double]"1"
$ws.Range("G16").Value = [double]"21.16964733333333"
$ws.Range("H16").Value = [double]"63.508942"
$ws.Range("I16").Value = [double]"0.1866181294337677"
$ws.Range("J16").Value = [double]"0.1866181294337677"
$ws.Range("K16").Value = [double]"3"
$ws.Range("L16").Value = [double]"1"
$ws.Range("M16").Value = [double]"61.37607633333332"
$ws.Range("N16").Value = [double]"184.128229"
$ws.Range("O16").Value = [double]"0.8280715929891429"
$ws.Range("P16").Value = [double]"0.8280715929891428"
$ws.Range("Q16").Value = [double]"1299.309890680413"
$ws.Range("R16").Value = [double]"11693.78901612372"
$ws.Range("S16").Value = [double]"0.1545331717208741"
$ws.Range("T16").Value = [double]"0.1545331717208741"

# Row 17
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Sdc2"
$ws.Range("C17").Value = "Ptprj"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = [double]"3"
$ws.Range("F17").Value = [double]"1"
$ws.Range("G17").Value = [double]"21.16964733333333"
$ws.Range("H17").Value = [double]"63.508942"
$ws.Range("I17").Value = [double]"0.1866181294337677"
$ws.Range("J17").Value = [double]"0.1866181294337677"
$ws.Range("K17").Value = [double]"3"
$ws.Range("L17").Value = [double]"1"
$ws.Range("M17").Value = [double]"5.243698"
$ws.Range("N17").Value = [double]"15.731094"
$ws.Range("O17").Value = [double]"0.07074674067517345"
$ws.Range("P17").Value = [double]"0.07074674067517343"
$ws.Range("Q17").Value = [double]"111.0072373825053"
$ws.Range("R17").Value = [double]"999.065136442548"
$ws.Range("S17").Value = [double]"0.01320262440833672"
$ws.Range("T17").Value = [double]"0.01320262440833671"
